# Apply the edits described by the commit diff:
#  - Follow limit per each go (B2): 10 -> 5
#  - Time limit to wait until next go (C2): 60 -> 5
#  - Active selection moves from B3 to C3
#  - Cosmetic workbook window position (xWindow) shifts 7680 -> 8610

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the parameter values on row 2 of the table.
$ws.Range("B2").Value = 5
$ws.Range("C2").Value = 5

# Move the active cell / selection to C3 (matches the saved <selection>).
$ws.Range("C3").Select() | Out-Null

# Best-effort nudge of the workbook window's on-screen position so it
# matches the recorded xWindow, in case the host honors it.
$excel.ActiveWindow.Left = 8610
$excel.ActiveWindow.Top = 105
